# Compose Mail: add localization strings for "Send School" / "Send Monitor"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 116: new key for the "send to school" compose-mail action
$ws.Range("A116").Value = "lang_mail_send_school"
$ws.Range("B116").Value = "Gửi Nhà Trường"
$ws.Range("C116").Value = "To School"

# Row 117: new key reusing the existing "Pick Up" translations
$ws.Range("A117").Value = "lang_pick_up"
$ws.Range("B117").Value = "Tuyến Đón"
$ws.Range("C117").Value = "Pick Up"

# Row 118: new key reusing the existing "Drop Down" translations
$ws.Range("A118").Value = "lang_drop_down"
$ws.Range("B118").Value = "Tuyến Trả"
$ws.Range("C118").Value = "Drop Down"

# Match the workbook's resulting view/selection state
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C118").Select()
